$wb = $excel.ActiveWorkbook

# Sheet "建物" (building): property_category column (I) rows 2-7 were mistakenly
# tagged as "land"; fix them to read "building".
$wsBuilding = $wb.Worksheets.Item("建物")
for ($r = 2; $r -le 7; $r++) {
    $wsBuilding.Cells.Item($r, 9).Value = "building"
}

# Sheet "汽車" (car): property_category column (H) row 2 was mistakenly tagged
# as "land"; fix it to read "car".
$wsCar = $wb.Worksheets.Item("汽車")
$wsCar.Cells.Item(2, 8).Value = "car"
